# Apply the 06/2022 lab-variable extraction update to selected_variables.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final (A,B) content for rows 1..61, taken from the target workbook state
$data = @(
  ,('included', 'excluded')
  ,('sodium', 'Activite anti-Xa (DOAC)')
  ,('chlore', 'activite anti-Xa (HBPM), therapeutique, 2x /jour')
  ,('hemoglobine glyquee', 'activite anti-Xa (HNF)')
  ,('ALAT', 'G-Sgv-pH(T), ABL')
  ,('calcium corrige', 'G-Sga-pO2(T), ABL')
  ,('triglycerides', 'G-Sga-pO2(a,T)/FO2(I), ABL')
  ,('fibrinogene', 'G-Sga-pCO2(T), ABL')
  ,('INR', 'pain scale')
  ,('proBNP', 'glycemie moyenne estimee')
  ,('PTT', 'cholesterol total')
  ,('potassium', 'osmolality')
  ,('bilirubine totale', 'G-Sga-calcium ionise (7.4), ABL')
  ,('creatinine', 'G-Sga-HCO3-(P,st), ABL')
  ,('cholesterol HDL', 'G-Sgv-HCO3-(P,st), ABL')
  ,('LDL cholesterol calcule', 'G-Sga-cBase(Ecf), ABL')
  ,('proteine C-reactive', 'G-Sgv-cBase(Ecf), ABL')
  ,('hemoglobine', 'cholesterol total')
  ,('glucose', 'G-Sgv-FCOHb, ABL')
  ,('erythrocytes', 'G-Sga-FCOHb, ABL')
  ,('uree', 'G-Sga-pH(T), ABL')
  ,('leucocytes', '')
  ,('ASAT', '')
  ,('hematocrite', '')
  ,('thrombocytes', '')
  ,('lactate', '')
  ,('cholesterol HDL', '')
  ,('NIHSS', '')
  ,('Glasgow Coma Scale', '')
  ,('FIO2', '')
  ,('oxygen_saturation', '')
  ,('systolic_blood_pressure', '')
  ,('diastolic_blood_pressure', '')
  ,('mean_blood_pressure', '')
  ,('heart_rate', '')
  ,('respiratory_rate', '')
  ,('temperature', '')
  ,('weight', '')
  ,('age', '')
  ,('Sex', '')
  ,('Referral', '')
  ,('Prestroke disability (Rankin)', '')
  ,('NIH on admission', '')
  ,('Antihypert. drugs pre-stroke', '')
  ,('Lipid lowering drugs pre-stroke', '')
  ,('Antiplatelet drugs', '')
  ,('Anticoagulants', '')
  ,('MedHist Hypertension', '')
  ,('MedHist Diabetes', '')
  ,('MedHist Hyperlipidemia', '')
  ,('MedHist Smoking', '')
  ,('MedHist Atrial Fibr.', '')
  ,('MedHist CHD', '')
  ,('MedHist PAD', '')
  ,('1st cholesterol total', '')
  ,('1st cholesterol LDL', '')
  ,('MedHist cerebrovascular_event', '')
  ,('categorical_onset_to_admission_time', '')
  ,('wake_up_stroke', '')
  ,('categorical_IVT', '')
  ,('categorical_IAT', '')
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 1
  $a = $data[$i][0]
  $b = $data[$i][1]
  $ws.Cells.Item($r, 1).Value = $a
  if ($b -ne '') {
    $ws.Cells.Item($r, 2).Value = $b
  }
}

# Reset view: clear the old scroll position / selection, select F19 (matches new sheetView)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$null = $ws.Range("F19").Select()

# Match the new print/page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "applied"
